$wb = $excel.ActiveWorkbook

# --- Instructions sheet: add two more bullet rows ---
$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsInstructions.Unprotect()
$wsInstructions.Range("A9").Value = "- Antibody details: Measurements or characteristics of the antibody"
$wsInstructions.Range("A10").Value = "- Antibody comment: Other comments on the antibody"
$wsInstructions.Protect()

# --- Antibodies sheet: add two new header columns ---
$wsAntibodies = $wb.Worksheets.Item("Antibodies")

# Copy format from the existing "Isotype" column (C) so the new header
# cells pick up the same bold header style, and so the sheet's used range
# naturally extends to row 2 (matching the existing C2) and column E.
$wsAntibodies.Range("C1:C2").Copy()
$wsAntibodies.Range("D1").PasteSpecial(-4122)
$wsAntibodies.Range("C1:C2").Copy()
$wsAntibodies.Range("E1").PasteSpecial(-4122)

$wsAntibodies.Range("D1").Value = "Antibody details"
$wsAntibodies.Range("E1").Value = "Antibody comment"

# New columns are a bit wider than the originals.
$wsAntibodies.Columns.Item(4).ColumnWidth = 15.166666666666666
$wsAntibodies.Columns.Item(5).ColumnWidth = 15.166666666666666
